# Insert a new data row at row 112 (pushing existing rows 112..217 down to 113..218)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(112).Insert()

$ws.Range("A112").Value2 = 1
$ws.Range("B112").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C112").Value2 = "Arica y Parinacota"
$ws.Range("D112").Value2 = 44658
$ws.Range("E112").Value2 = 15
$ws.Range("F112").Value2 = "Fruta"
$ws.Range("G112").Value2 = 100108
$ws.Range("H112").Value2 = "Tropicales y subtropicales"
$ws.Range("I112").Value2 = 100108006
$ws.Range("J112").Value2 = "Plátano"
$ws.Range("K112").Value2 = "Sin especificar"
$ws.Range("L112").Value2 = "Pintón"
$ws.Range("M112").Value2 = 120
$ws.Range("N112").Value2 = 17000
$ws.Range("O112").Value2 = 18000
$ws.Range("P112").Value2 = 17500
$ws.Range("Q112").Value2 = "$/caja 20 kilos"
$ws.Range("R112").Value2 = "Ecuador"
$ws.Range("S112").Value2 = 875
$ws.Range("T112").Value2 = 20
